$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new data row (row 11), mirroring the format of the prior row (row 10)
# Force column A to be stored as text so the date-like string is not
# auto-converted into a date serial number.
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "02-10-2025"
$ws.Range("B11").Value = "The price of gold in India today is ₹11,869 per gram for 24 karat gold, ₹10,880 per gram for 22 karat gold and ₹8,902 per gram for 18 karat gold (also called 999 gold)."

# Copy the formatting (borders, wrap text, etc.) from row 10 down to row 11
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to match the new active cell, same as the source workbook
$ws.Range("B11").Select()
